# Update res_bus/vm_pu.xlsx values for Case_5_10 (case with 380 kV done)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.065372528896858
$ws.Range("D2").Value = 1.063869726135687
$ws.Range("E2").Value = 1.06945816529658
$ws.Range("F2").Value = 1.077626730693453
$ws.Range("I2").Value = 1.042356725517955
$ws.Range("J2").Value = 1.070327980816843
$ws.Range("K2").Value = 1.066587198548006
$ws.Range("L2").Value = 1.072160619036924
$ws.Range("M2").Value = 1.08030753462386
$ws.Range("N2").Value = 1.005712725503983
# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.066880291374384
$ws.Range("D3").Value = 1.065008334230027
$ws.Range("E3").Value = 1.070778536756293
$ws.Range("F3").Value = 1.078978983542374
$ws.Range("I3").Value = 1.042658868434196
$ws.Range("J3").Value = 1.071488678508449
$ws.Range("K3").Value = 1.067540090210302
$ws.Range("L3").Value = 1.073295914759028
$ws.Range("M3").Value = 1.081476212547096
# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.067854960537445
$ws.Range("D4").Value = 1.065744020889579
$ws.Range("E4").Value = 1.071632255198699
$ws.Range("F4").Value = 1.079853312558215
$ws.Range("I4").Value = 1.042852545510387
$ws.Range("J4").Value = 1.072238341950822
$ws.Range("K4").Value = 1.068155022228101
$ws.Range("L4").Value = 1.074029328100156
$ws.Range("M4").Value = 1.082231211260302
# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.068264489347675
$ws.Range("D5").Value = 1.066053051444677
$ws.Range("E5").Value = 1.071991007378773
$ws.Range("F5").Value = 1.080220724962723
$ws.Range("I5").Value = 1.042933530659425
$ws.Range("J5").Value = 1.072553173250725
$ws.Range("K5").Value = 1.068413147586797
$ws.Range("L5").Value = 1.074337372281106
$ws.Range("M5").Value = 1.082548326332788
# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.068333238194525
$ws.Range("D6").Value = 1.066104924360007
$ws.Range("E6").Value = 1.072051234723898
$ws.Range("F6").Value = 1.080282406137341
$ws.Range("I6").Value = 1.042947102852417
$ws.Range("J6").Value = 1.072606015753812
$ws.Range("K6").Value = 1.068456465083194
$ws.Range("L6").Value = 1.074389077794368
$ws.Range("M6").Value = 1.082601554678379
# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.067860433546391
$ws.Range("D7").Value = 1.065748151156951
$ws.Range("E7").Value = 1.071637049450723
$ws.Range("F7").Value = 1.079858222546054
$ws.Range("I7").Value = 1.042853629352201
$ws.Range("J7").Value = 1.072242550022213
$ws.Range("K7").Value = 1.068158472847897
$ws.Range("L7").Value = 1.074033445307505
$ws.Range("M7").Value = 1.082235449687457
# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.065882285600336
$ws.Range("D8").Value = 1.064254747681384
$ws.Range("E8").Value = 1.069904528414337
$ws.Range("F8").Value = 1.078083872198963
$ws.Range("I8").Value = 1.042459215754636
$ws.Range("J8").Value = 1.070720533643294
$ws.Range("K8").Value = 1.06690957708747
$ws.Range("L8").Value = 1.072544548336917
$ws.Range("M8").Value = 1.080702748746967
# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.062388947969976
$ws.Range("D9").Value = 1.061614823044059
$ws.Range("E9").Value = 1.066846418508551
$ws.Range("F9").Value = 1.074951907571959
$ws.Range("I9").Value = 1.041750136788557
$ws.Range("J9").Value = 1.068027732235105
$ws.Range("K9").Value = 1.064696045314727
$ws.Range("L9").Value = 1.069911549159657
$ws.Range("M9").Value = 1.077992444153965
# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.060054558507437
$ws.Range("D10").Value = 1.059849014048127
$ws.Range("E10").Value = 1.064803900057788
$ws.Range("F10").Value = 1.072860051467305
$ws.Range("I10").Value = 1.041267872235626
$ws.Range("J10").Value = 1.066224975744576
$ws.Range("K10").Value = 1.063211512394594
$ws.Range("L10").Value = 1.068149646741638
$ws.Range("M10").Value = 1.076178924115094
# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.059042349516443
$ws.Range("D11").Value = 1.059082955269733
$ws.Range("E11").Value = 1.063918505300229
$ws.Range("F11").Value = 1.071953268280647
$ws.Range("I11").Value = 1.041056763240856
$ws.Range("J11").Value = 1.065442506219419
$ws.Range("K11").Value = 1.062566544148675
$ws.Range("L11").Value = 1.067385105036827
$ws.Range("M11").Value = 1.075392011819893
# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.058666151541673
$ws.Range("D12").Value = 1.058798184041611
$ws.Range("E12").Value = 1.063589478687138
$ws.Range("F12").Value = 1.071616293249206
$ws.Range("I12").Value = 1.040978002916815
$ws.Range("J12").Value = 1.065151576543332
$ws.Range("K12").Value = 1.062326645857333
$ws.Range("L12").Value = 1.067100870590839
$ws.Range("M12").Value = 1.075099464919254
# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.058746857295588
$ws.Range("D13").Value = 1.058859278599849
$ws.Range("E13").Value = 1.063660062964155
$ws.Range("F13").Value = 1.071688582662069
$ws.Range("I13").Value = 1.040994912915216
$ws.Range("J13").Value = 1.065213994963237
$ws.Range("K13").Value = 1.062378119787078
$ws.Range("L13").Value = 1.067161851225113
$ws.Range("M13").Value = 1.075162228740698
# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.059011257376327
$ws.Range("D14").Value = 1.059059420565756
$ws.Range("E14").Value = 1.063891310983165
$ws.Range("F14").Value = 1.071925417022797
$ws.Range("I14").Value = 1.041050259936621
$ws.Range("J14").Value = 1.065418463741439
$ws.Range("K14").Value = 1.062546720805813
$ws.Range("L14").Value = 1.067361615251593
$ws.Range("M14").Value = 1.075367834986455
# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.059174133809283
$ws.Range("D15").Value = 1.059182704958274
$ws.Range("E15").Value = 1.064033770262268
$ws.Range("F15").Value = 1.072071317725526
$ws.Range("I15").Value = 1.041084315290195
$ws.Range("J15").Value = 1.065544405683792
$ws.Range("K15").Value = 1.06265055781097
$ws.Range("L15").Value = 1.067484663226367
$ws.Range("M15").Value = 1.075494482143669
# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.060121705269047
$ws.Range("D16").Value = 1.059899823851724
$ws.Range("E16").Value = 1.064862639850438
$ws.Range("F16").Value = 1.072920210209623
$ws.Range("I16").Value = 1.041281834550928
$ws.Range("J16").Value = 1.066276865919066
$ws.Range("K16").Value = 1.06325427101164
$ws.Range("L16").Value = 1.068200352142035
$ws.Range("M16").Value = 1.07623111375661
# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.060715710773627
$ws.Range("D17").Value = 1.060349261535748
$ws.Range("E17").Value = 1.065382303904695
$ws.Range("F17").Value = 1.073452427387219
$ws.Range("I17").Value = 1.04140512013282
$ws.Range("J17").Value = 1.066735815867276
$ws.Range("K17").Value = 1.063632384045326
$ws.Range("L17").Value = 1.068648845760389
$ws.Range("M17").Value = 1.076692738838733
# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.061062049294343
$ws.Range("D18").Value = 1.06061127117706
$ws.Range("E18").Value = 1.06568532187809
$ws.Range("F18").Value = 1.073762764986949
$ws.Range("I18").Value = 1.041476810148867
$ws.Range("J18").Value = 1.067003334132587
$ws.Range("K18").Value = 1.063852723293632
$ws.Range("L18").Value = 1.06891028794177
$ws.Range("M18").Value = 1.076961838471456
# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.061180119140652
$ws.Range("D19").Value = 1.06070058613215
$ws.Range("E19").Value = 1.065788627527203
$ws.Range("F19").Value = 1.073868566034592
$ws.Range("I19").Value = 1.041501217263119
$ws.Range("J19").Value = 1.067094520674274
$ws.Range("K19").Value = 1.06392781819121
$ws.Range("L19").Value = 1.068999406573157
$ws.Range("M19").Value = 1.077053567773304
# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.060651993607726
$ws.Range("D20").Value = 1.060301055615953
$ws.Range("E20").Value = 1.065326558559462
$ws.Range("F20").Value = 1.073395335459162
$ws.Range("I20").Value = 1.041391915563812
$ws.Range("J20").Value = 1.06668659348587
$ws.Range("K20").Value = 1.063591837609975
$ws.Range("L20").Value = 1.06860074287154
$ws.Range("M20").Value = 1.076643227265906
# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.058933404236012
$ws.Range("D21").Value = 1.059000489936497
$ws.Range("E21").Value = 1.063823218435757
$ws.Range("F21").Value = 1.071855679542489
$ws.Range("I21").Value = 1.041033971156992
$ws.Range("J21").Value = 1.065358260713143
$ws.Range("K21").Value = 1.062497081108906
$ws.Range("L21").Value = 1.067302796657957
$ws.Range("M21").Value = 1.075307296088613
# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.057851590935971
$ws.Range("D22").Value = 1.058181480720051
$ws.Range("E22").Value = 1.062877128984111
$ws.Range("F22").Value = 1.070886734873252
$ws.Range("I22").Value = 1.040806920128907
$ws.Range("J22").Value = 1.064521429664941
$ws.Range("K22").Value = 1.061806862160727
$ws.Range("L22").Value = 1.066485279519103
$ws.Range("M22").Value = 1.074465878792818
# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.058425203407647
$ws.Range("D23").Value = 1.058615777064443
$ws.Range("E23").Value = 1.063378754224119
$ws.Range("F23").Value = 1.071400478206472
$ws.Range("I23").Value = 1.040927474057438
$ws.Range("J23").Value = 1.06496520843903
$ws.Range("K23").Value = 1.062172942089881
$ws.Range("L23").Value = 1.066918799736613
$ws.Range("M23").Value = 1.074912070609507
# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.060680785069293
$ws.Range("D24").Value = 1.060322838228285
$ws.Range("E24").Value = 1.065351747770016
$ws.Range("F24").Value = 1.073421133139956
$ws.Range("I24").Value = 1.041397882821241
$ws.Range("J24").Value = 1.066708835516142
$ws.Range("K24").Value = 1.063610159440828
$ws.Range("L24").Value = 1.068622478980506
$ws.Range("M24").Value = 1.076665599904442
# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.063292999338938
$ws.Range("D25").Value = 1.062298321835474
$ws.Range("E25").Value = 1.067637657734827
$ws.Range("F25").Value = 1.075762257629973
$ws.Range("I25").Value = 1.041935126903194
$ws.Range("J25").Value = 1.068725196677536
$ws.Range("K25").Value = 1.065269838600203
$ws.Range("L25").Value = 1.07059338074378
$ws.Range("M25").Value = 1.078694275121934
